# Add a new column J (year 2022) to the primary-education completion rate
# table on Sheet1, mirroring the formatting of column I (2021) and filling
# in the 2022 values for each region, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column I (rows 3-14, the data block incl. the
# blank/border row 3) across to column J so the new column visually
# matches the rest of the table (fonts, borders, number formats, etc).
$ws.Range("I3:I14").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header year
$ws.Range("J4").Value = 2022

# 2022 data values for each row (Kyrgyz Republic + oblasts + cities)
$ws.Range("J5").Value = 96.4
$ws.Range("J6").Value = 96.4
$ws.Range("J7").Value = 97.9
$ws.Range("J8").Value = 95.3
$ws.Range("J9").Value = 93.8
$ws.Range("J10").Value = 95.5
$ws.Range("J11").Value = 94.4
$ws.Range("J12").Value = 95
$ws.Range("J13").Value = 98.7
$ws.Range("J14").Value = 97.3

# Match the saved cursor/selection position recorded in the workbook.
$ws.Range("L10").Select() | Out-Null
